$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the longitude/latitude header cells to capitalized versions
$ws.Range("C1").Value = "Longitude"
$ws.Range("D1").Value = "Latitude"

$ws.Range("D1").Select()
